$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3642143333333334
$ws.Range("H2").Value = 1.092643
$ws.Range("I2").Value = 0.4800482050304226
$ws.Range("J2").Value = 0.4800482050304224
$ws.Range("M2").Value = 30.46625333333334
$ws.Range("N2").Value = 91.39876000000001
$ws.Range("O2").Value = 0.2185380492512374
$ws.Range("P2").Value = 0.2331534018544084
$ws.Range("Q2").Value = 11.09624614696445
$ws.Range("R2").Value = 99.86621532268002
$ws.Range("S2").Value = 0.1049087982739066
$ws.Range("T2").Value = 0.1119248720569455

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3642143333333334
$ws.Range("H3").Value = 1.092643
$ws.Range("I3").Value = 0.4800482050304226
$ws.Range("J3").Value = 0.4800482050304224
$ws.Range("O3").Value = 0.2491807703757967
$ws.Range("P3").Value = 0.2658454419670822
$ws.Range("Q3").Value = 12.65212704448267
$ws.Range("R3").Value = 113.869143400344
$ws.Range("S3").Value = 0.1196187815469991
$ws.Range("T3").Value = 0.1276186272318172

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3642143333333334
$ws.Range("H4").Value = 1.092643
$ws.Range("I4").Value = 0.4800482050304226
$ws.Range("J4").Value = 0.4800482050304224
$ws.Range("M4").Value = 23.69037333333334
$ws.Range("N4").Value = 71.07112000000001
$ws.Range("O4").Value = 0.1699338582153697
$ws.Range("P4").Value = 0.181298667526812
$ws.Range("Q4").Value = 8.62837353001778
$ws.Range("R4").Value = 77.65536177016001
$ws.Range("S4").Value = 0.08157644361018257
$ws.Range("T4").Value = 0.08703209992065344

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3642143333333334
$ws.Range("H5").Value = 1.092643
$ws.Range("I5").Value = 0.4800482050304226
$ws.Range("J5").Value = 0.4800482050304224
$ws.Range("M5").Value = 26.2168665
$ws.Range("N5").Value = 52.433733
$ws.Range("O5").Value = 0.18805669340777
$ws.Range("P5").Value = 0.1337556791894743
$ws.Range("Q5").Value = 9.548558554386501
$ws.Range("R5").Value = 57.291351326319
$ws.Range("S5").Value = 0.09027627811435648
$ws.Range("T5").Value = 0.06420917370753218

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.3642143333333334
$ws.Range("H6").Value = 1.092643
$ws.Range("I6").Value = 0.4800482050304226
$ws.Range("J6").Value = 0.4800482050304224
$ws.Range("M6").Value = 24.297748
$ws.Range("N6").Value = 72.893244
$ws.Range("O6").Value = 0.1742906287498262
$ws.Range("P6").Value = 0.1859468094622229
$ws.Range("Q6").Value = 8.849588089321333
$ws.Range("R6").Value = 79.646292803892
$ws.Range("S6").Value = 0.08366790348497782
$ws.Range("T6").Value = 0.0892634321134741

# Row 7
$ws.Range("G7").Value = 0.3944893333333333
$ws.Range("H7").Value = 1.183468
$ws.Range("I7").Value = 0.5199517949695774
$ws.Range("J7").Value = 0.5199517949695774
$ws.Range("M7").Value = 30.46625333333334
$ws.Range("N7").Value = 91.39876000000001
$ws.Range("O7").Value = 0.2185380492512374
$ws.Range("P7").Value = 0.2331534018544084
$ws.Range("Q7").Value = 12.01861196663111
$ws.Range("R7").Value = 108.16750769968
$ws.Range("S7").Value = 0.1136292509773308
$ws.Range("T7").Value = 0.1212285297974629

# Row 8
$ws.Range("G8").Value = 0.3944893333333333
$ws.Range("H8").Value = 1.183468
$ws.Range("I8").Value = 0.5199517949695774
$ws.Range("J8").Value = 0.5199517949695774
$ws.Range("O8").Value = 0.2491807703757967
$ws.Range("P8").Value = 0.2658454419670822
$ws.Range("Q8").Value = 13.70382411188266
$ws.Range("R8").Value = 123.334417006944
$ws.Range("S8").Value = 0.1295619888287976
$ws.Range("T8").Value = 0.138226814735265

# Row 9
$ws.Range("G9").Value = 0.3944893333333333
$ws.Range("H9").Value = 1.183468
$ws.Range("I9").Value = 0.5199517949695774
$ws.Range("J9").Value = 0.5199517949695774
$ws.Range("M9").Value = 23.69037333333334
$ws.Range("N9").Value = 71.07112000000001
$ws.Range("O9").Value = 0.1699338582153697
$ws.Range("P9").Value = 0.181298667526812
$ws.Range("Q9").Value = 9.345599582684445
$ws.Range("R9").Value = 84.11039624416
$ws.Range("S9").Value = 0.08835741460518717
$ws.Range("T9").Value = 0.09426656760615855

# Row 10
$ws.Range("G10").Value = 0.3944893333333333
$ws.Range("H10").Value = 1.183468
$ws.Range("I10").Value = 0.5199517949695774
$ws.Range("J10").Value = 0.5199517949695774
$ws.Range("M10").Value = 26.2168665
$ws.Range("N10").Value = 52.433733
$ws.Range("O10").Value = 0.18805669340777
$ws.Range("P10").Value = 0.1337556791894743
$ws.Range("Q10").Value = 10.342274187674
$ws.Range("R10").Value = 62.053645126044
$ws.Range("S10").Value = 0.0977804152934135
$ws.Range("T10").Value = 0.06954650548194212

# Row 11
$ws.Range("G11").Value = 0.3944893333333333
$ws.Range("H11").Value = 1.183468
$ws.Range("I11").Value = 0.5199517949695774
$ws.Range("J11").Value = 0.5199517949695774
$ws.Range("M11").Value = 24.297748
$ws.Range("N11").Value = 72.893244
$ws.Range("O11").Value = 0.1742906287498262
$ws.Range("P11").Value = 0.1859468094622229
$ws.Range("Q11").Value = 9.585202410021331
$ws.Range("R11").Value = 86.26682169019199
$ws.Range("S11").Value = 0.09062272526484838
$ws.Range("T11").Value = 0.09668337734874882
